$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G7").ClearContents()
Write-Output "done"
